# Updates cryptos list prices (D) and 1h volume deltas (E) per latest scrape.
# D-column values are numeric-looking text (e.g. "251.10", "37.254.26"); a
# leading apostrophe forces them to stay text (matching the source data, which
# stores prices as literal strings) instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.254.26"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "'2.096.67"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'251.10"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'0.661"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'52.51"
$ws.Range("E8").Value = "  +15.31%  "
$ws.Range("D9").Value = "'61.56"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").Value = "'0.375"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "'0.0743"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  +7.25%  "
$ws.Range("D13").Value = "'15.14"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "'2.399.36"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "'0.834"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "'2.092.78"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "'5.14"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "'37.222.71"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'72.41"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'14.21"
$ws.Range("E20").Value = "  +8.94%  "
$ws.Range("D21").Value = "'0.0₃0840"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'240.96"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "'5.24"
$ws.Range("E23").Value = "  +6.50%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "'170.61"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  +6.54%  "
$ws.Range("D28").Value = "'20.73"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "'2.00"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'22.23"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  +25.50%  "
$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'0.0611"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "'0.0906"
$ws.Range("E35").Value = "  +10.10%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'1.86"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E38").Value = "  +6.05%  "
$ws.Range("D39").Value = "'4.11"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "'1.35"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +11.83%  "
$ws.Range("D42").Value = "'0.0225"
$ws.Range("E42").Value = "  +3.18%  "
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("D44").Value = "'98.80"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "'0.0921"
$ws.Range("E45").Value = "  +13.00%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'1.321.32"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'2.97"
$ws.Range("E48").Value = "  +6.65%  "
$ws.Range("D49").Value = "'3.88"
$ws.Range("E49").Value = "  +83.69%  "
$ws.Range("D50").Value = "'7.03"
$ws.Range("E50").Value = "  +13.58%  "
$ws.Range("D51").Value = "'2.284.94"
$ws.Range("E51").Value = "  +1.23%  "
